$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank column before column A; every existing column
# (and the three merged header cells) shifts one slot to the right.
$ws.Range("A1").EntireColumn.Insert() | Out-Null

# The old title "Geometry and Permittivity Definition" merge (originally
# A1:M1) is now B1:N1 after the shift. Detach it, move its text into the
# new A1, and re-merge across A1:N1 so the title spans the new column too.
$ws.Range("B1:N1").UnMerge() | Out-Null
$ws.Range("B1").Cut($ws.Range("A1")) | Out-Null
$ws.Range("A1:N1").Merge() | Out-Null

# New header cell for the inserted column: same look as the other
# "Geometry and Permittivity Definition" sub-headers in row 2 (B2:N2).
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$ws.Range("A2").Value = "Inclusion"

# Populate the new "Inclusion" id column for each data row.
$inclusionIds = @(223, 232, 332, 323, 442, 424, 224, 242, 225, 252)
for ($i = 0; $i -lt $inclusionIds.Length; $i++) {
    $ws.Cells.Item(3 + $i, 1).Value = $inclusionIds[$i]
}

$excel.CutCopyMode = 0

# Match the saved selection state.
$ws.Range("C7").Select() | Out-Null
